# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker "IVAN JOSE JIMENEZ RODRIGUEZ" (CC 1128045090) now has account
# statement rows for the earlier periods 1809-1901 as well (previously the
# sheet only covered 1901-2002 for him). Adding these 4 extra periods shifts
# the table so that the single existing rows for ENUAR YAMITH CASTILLO
# RAMIREZ / CRISTIAN ALBERTO RIVAS GALVAN (both period 1902) land further
# down, and the remaining IVAN JOSE rows (1902-2002) follow them, all still
# sorted by Periodo Mora ascending.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Tipo Doc, N Doc Trabajador, Nombre Trabajador, Periodo Mora, Valor Mora, Salario Basico
$data = @(
    @("CC", "1128045090", "IVAN JOSE JIMENEZ RODRIGUEZ",   "1809", 31249, 781242),
    @("CC", "1128045090", "IVAN JOSE JIMENEZ RODRIGUEZ",   "1810", 31249, 781242),
    @("CC", "1128045090", "IVAN JOSE JIMENEZ RODRIGUEZ",   "1811", 31249, 781242),
    @("CC", "1128045090", "IVAN JOSE JIMENEZ RODRIGUEZ",   "1812", 31249, 781242),
    @("CC", "1128045090", "IVAN JOSE JIMENEZ RODRIGUEZ",   "1901", 31249, 781242),
    @("CC", "73434861",   "ENUAR YAMITH CASTILLO RAMIREZ", "1902", 31249, 828116),
    @("CC", "1047405914", "CRISTIAN ALBERTO RIVAS GALVAN", "1902", 31249, 781242),
    @("CC", "1128045090", "IVAN JOSE JIMENEZ RODRIGUEZ",   "1902", 31249, 781242),
    @("CC", "1128045090", "IVAN JOSE JIMENEZ RODRIGUEZ",   "1903", 31249, 781242),
    @("CC", "1128045090", "IVAN JOSE JIMENEZ RODRIGUEZ",   "1904", 31249, 781242),
    @("CC", "1128045090", "IVAN JOSE JIMENEZ RODRIGUEZ",   "1905", 31249, 781242),
    @("CC", "1128045090", "IVAN JOSE JIMENEZ RODRIGUEZ",   "1906", 31249, 781242),
    @("CC", "1128045090", "IVAN JOSE JIMENEZ RODRIGUEZ",   "1907", 31249, 781242),
    @("CC", "1128045090", "IVAN JOSE JIMENEZ RODRIGUEZ",   "1908", 31249, 781242),
    @("CC", "1128045090", "IVAN JOSE JIMENEZ RODRIGUEZ",   "1909", 31249, 781242),
    @("CC", "1128045090", "IVAN JOSE JIMENEZ RODRIGUEZ",   "1910", 31249, 781242),
    @("CC", "1128045090", "IVAN JOSE JIMENEZ RODRIGUEZ",   "1911", 31249, 781242),
    @("CC", "1128045090", "IVAN JOSE JIMENEZ RODRIGUEZ",   "1912", 31249, 781242),
    @("CC", "1128045090", "IVAN JOSE JIMENEZ RODRIGUEZ",   "2001", 31249, 781242),
    @("CC", "1128045090", "IVAN JOSE JIMENEZ RODRIGUEZ",   "2002", 17708, 781242)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]
    $ws.Range("G$row").Value = $vals[5]
}
